$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new product row for "FLOPADEX 8 MG 30 CAPSULES" right after
#    "FLACORT 30MG 20 TAB" (currently row 26), i.e. at row 27, pushing the
#    rest of the table (starting with "FLUREST N 20 TABS") down by one.
# ---------------------------------------------------------------------------
$ws.Rows("27:27").Insert(-4121)

# Clone the formatting (styles/borders/fill) of the row directly above so the
# new row matches the existing table look exactly.
$ws.Range("A26:N26").Copy()
$ws.Range("A27:N27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Rows("27:27").RowHeight = 25.5

$ws.Range("A27").Value = 24
$ws.Range("B27").Value = "FLOPADEX 8 MG 30 CAPSULES"
$ws.Range("H27").Value = "1:1"
$ws.Range("L27").Value = 59
$ws.Range("N27").Value = "0:0"

# ---------------------------------------------------------------------------
# 2) Insert a new product row for "TARGOFLOXIN 750 MG 10 F.C.TABS." right
#    after "T4-THYRO 100MCG 100 TABLETS" (now row 57, having shifted down by
#    one because of the previous insert), i.e. at row 58, pushing the rest of
#    the table (starting with "TAVONIZA 20 MG 20 F.C.TABS.") down by one.
# ---------------------------------------------------------------------------
$ws.Rows("58:58").Insert(-4121)

$ws.Range("A57:N57").Copy()
$ws.Range("A58:N58").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Rows("58:58").RowHeight = 25.5

$ws.Range("A58").Value = 55
$ws.Range("B58").Value = "TARGOFLOXIN 750 MG 10 F.C.TABS."
$ws.Range("H58").Value = "1:0"
$ws.Range("L58").Value = 92.12
$ws.Range("N58").Value = "0:2"

$wb.Save()
